$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Project / Release" column J with the next version bump for each
# Tardigrade package. Cell values are written in the same order the
# original author entered them (this controls shared-string insertion
# order), then formatting (bold header / italic minor-bump cells) is
# applied afterwards.
$ws.Range("J1").Value = "11.5.0"
$ws.Range("J2").Value = "9.1.1"
$ws.Range("J5").Value = "2.1.0"
$ws.Range("J3").Value = "5.1.1"
$ws.Range("J8").Value = "8.3.2"
$ws.Range("J4").Value = "3.4.1"
$ws.Range("J9").Value = "1.3.2"
$ws.Range("J10").Value = "3.3.2"
$ws.Range("J11").Value = "1.1.2"
$ws.Range("J6").Value = "6.1.2"
$ws.Range("J7").Value = "10.0.1"

# Header (J1) uses the same bold style as the rest of row 1.
$ws.Range("J1").Font.Bold = $true

# Some cells use the italic style reserved for notable bumps.
$ws.Range("J3").Font.Italic = $true
$ws.Range("J6").Font.Italic = $true
$ws.Range("J7").Font.Italic = $true

# Selection ends on J5, matching the saved view state.
$ws.Range("J5").Select() | Out-Null
